$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

# Remove the "Event" (M) and "Correction" (N) columns' data for rows 1-12.
# This clears the cell contents (not a column delete/shift) so the sheet's
# used range shrinks from A1:N12 down to A1:L12.
$ws.Range("M1:N12").Clear()
